# Insert a new data row at row 54, shifting existing rows 54-138 down to 55-139.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(54).Insert()

$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44580
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100103
$ws.Range("H54").Value = "Frutos de hueso (carozo)"
$ws.Range("I54").Value = 100103004
$ws.Range("J54").Value = "Durazno"
$ws.Range("K54").Value = "Carson"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 220
$ws.Range("N54").Value = 11000
$ws.Range("O54").Value = 12000
$ws.Range("P54").Value = 11545
$ws.Range("Q54").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R54").Value = "Región de O'Higgins"
$ws.Range("S54").Value = 722
$ws.Range("T54").Value = 16
